$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CONFIG")

# Row/col layout for the refreshed annotation table (A1:I5).
# Values are entered as ="literal" formulas and then converted to
# plain static text via copy / paste-special-values. Doing it this
# way (instead of Range.Value, which auto-detects "5500.0" etc. as
# numbers) keeps every cell a genuine text value - matching the
# t="s" shared-string cells in the target workbook - without leaving
# behind any NumberFormat/quote-prefix style residue.

$ws.Range("A1").Formula = "=""5.5 GHz"""
$ws.Range("B1").Formula = "=""5500.0"""
$ws.Range("C1").Formula = "=""1000.0"""
$ws.Range("D1").Formula = "=""5000.0"""
$ws.Range("E1").Formula = "=""6000.0"""
$ws.Range("F1").Formula = "=""0.1571"""
$ws.Range("G1").Formula = "=""10.0"""
$ws.Range("H1").Formula = "=""20"""
$ws.Range("I1").Formula = "=""-50"""

$ws.Range("A2").Formula = "=""4 GHz"""
$ws.Range("B2").Formula = "=""4000.0"""
$ws.Range("C2").Formula = "=""1000.0"""
$ws.Range("D2").Formula = "=""3500.0"""
$ws.Range("E2").Formula = "=""4500.0"""
$ws.Range("F2").Formula = "=""0.1571"""
$ws.Range("G2").Formula = "=""10.0"""
$ws.Range("H2").Formula = "=""20"""
$ws.Range("I2").Formula = "=""-50"""

$ws.Range("A3").Formula = "=""915 MHz"""
$ws.Range("B3").Formula = "=""915.0"""
$ws.Range("C3").Formula = "=""100.0"""
$ws.Range("D3").Formula = "=""865.0"""
$ws.Range("E3").Formula = "=""965.0"""
$ws.Range("F3").Formula = "=""0.03945"""
$ws.Range("G3").Formula = "=""10.0"""
$ws.Range("H3").Formula = "=""20"""
$ws.Range("I3").Formula = "=""-50"""

$ws.Range("A4").Formula = "=""863 MHz"""
$ws.Range("B4").Formula = "=""863.0"""
$ws.Range("C4").Formula = "=""100.0"""
$ws.Range("D4").Formula = "=""813.0"""
$ws.Range("E4").Formula = "=""913.0"""
$ws.Range("F4").Formula = "=""0.03945"""
$ws.Range("G4").Formula = "=""10.0"""
$ws.Range("H4").Formula = "=""20"""
$ws.Range("I4").Formula = "=""-50"""

$ws.Range("A5").Formula = "=""Wide Band"""
$ws.Range("B5").Formula = "=""3015.0"""
$ws.Range("C5").Formula = "=""5970.0"""
$ws.Range("D5").Formula = "=""30.0"""
$ws.Range("E5").Formula = "=""6000.0"""
$ws.Range("F5").Formula = "=""0.3156"""
$ws.Range("G5").Formula = "=""1.0"""
$ws.Range("H5").Formula = "=""20"""
$ws.Range("I5").Formula = "=""-50"""

$dataRange = $ws.Range("A1:I5")
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
